# Update the "Förändrad" date column (C2:C12) from 2023-09-16 (45185)
# to 2023-10-05 (45204), as part of an automatic update of files.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
